$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 350
$ws.Range("D10").Value = 850
$ws.Range("B25").Value = 300
$ws.Range("D25").Value = 500
$ws.Range("D29").Value = 561
$ws.Range("D31").Value = 155

$ws.Range("H17").Select()
